$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-like numeric strings (e.g. "0.650", "175.50") are not
# auto-converted to numbers by Excel, which would drop trailing zeros.
$textCells = @("D2","D3","D5","D6","D7","D9","D10","D12","D14","D16","D17","D18","D19","D21","D22","D23","D24","D25","D28","D29","D30","D31","D32","D33","D36","D37","D40","D41","D43","D44","D46","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "43.190.26"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "2.395.35"
$ws.Range("E3").Value = "  +5.15%  "
$ws.Range("E4").Value = "  -0.41%  "
$ws.Range("D5").Value = "332.62"
$ws.Range("E5").Value = "  +9.08%  "
$ws.Range("D6").Value = "105.25"
$ws.Range("D7").Value = "0.650"
$ws.Range("E7").Value = "  +2.86%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "0.649"
$ws.Range("E9").Value = "  +5.74%  "
$ws.Range("D10").Value = "41.99"
$ws.Range("E10").Value = "  -6.45%  "
$ws.Range("E11").Value = "  +1.46%  "
$ws.Range("D12").Value = "8.72"
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("D14").Value = "17.06"
$ws.Range("E14").Value = "  +10.72%  "
$ws.Range("E15").Value = "  +1.95%  "
$ws.Range("D16").Value = "2.756.91"
$ws.Range("E16").Value = "  +5.15%  "
$ws.Range("D17").Value = "2.394.24"
$ws.Range("E17").Value = "  +1.43%  "
$ws.Range("D18").Value = "43.145.96"
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("D19").Value = "7.69"
$ws.Range("E19").Value = "  +6.32%  "
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("D21").Value = "3.81"
$ws.Range("E21").Value = "  +7.43%  "
$ws.Range("D22").Value = "77.26"
$ws.Range("E22").Value = "  +3.16%  "
$ws.Range("D23").Value = "275.34"
$ws.Range("E23").Value = "  +7.90%  "
$ws.Range("D24").Value = "2.42"
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("D25").Value = "9.71"
$ws.Range("E25").Value = "  +7.45%  "
$ws.Range("E26").Value = "  +1.76%  "
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("D28").Value = "23.27"
$ws.Range("E28").Value = "  +5.12%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.19"
$ws.Range("E29").Value = "  -2.13%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "175.50"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("B31").Value = "WEMIXToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D31").Value = "3.16"
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").Value = "37.05"
$ws.Range("E32").Value = "  -2.97%  "
$ws.Range("D33").Value = "0.0936"
$ws.Range("E33").Value = "  +4.23%  "
$ws.Range("E34").Value = "  +4.89%  "
$ws.Range("E35").Value = "  +4.84%  "
$ws.Range("D36").Value = "4.89"
$ws.Range("E36").Value = "  -3.17%  "
$ws.Range("D37").Value = "4.10"
$ws.Range("E37").Value = "  -3.83%  "
$ws.Range("E38").Value = "  -3.74%  "
$ws.Range("E39").Value = "  +1.75%  "
$ws.Range("D40").Value = "2.82"
$ws.Range("E40").Value = "  +11.74%  "
$ws.Range("D41").Value = "1.57"
$ws.Range("E41").Value = "  +14.26%  "
$ws.Range("E42").Value = "  +0.89%  "
$ws.Range("D43").Value = "70.07"
$ws.Range("E43").Value = "  -3.80%  "
$ws.Range("D44").Value = "122.68"
$ws.Range("E44").Value = "  +14.80%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "91.72"
$ws.Range("E46").Value = "  +45.30%  "
$ws.Range("D47").Value = "12.33"
$ws.Range("E47").Value = "  -2.42%  "
$ws.Range("D48").Value = "5.57"
$ws.Range("E48").Value = "  -1.08%  "
$ws.Range("D49").Value = "9.31"
$ws.Range("E49").Value = "  +5.92%  "
$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D50").Value = "0.509"
$ws.Range("E50").Value = "  +16.23%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "1.32"
$ws.Range("E51").Value = "  +1.53%  "
